# Update products and translations
# - Remove the "unit_en" column (column K): it shifts the old type_ar/type_en
#   columns (L, M) left into K, L.
# - Repurpose the old "unit_ar" column (J) into a new numeric "discount" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the unit_en column entirely (K). This also shifts L->K and M->L,
# and prunes the now-unused "kilo"/"lite" shared strings automatically.
$ws.Columns("K").Delete()

# Turn the former unit_ar column (J) into the new "discount" column with
# numeric values.
$ws.Range("J1").Value = "discount"
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 2

# Give column I a bit more width, and reflect the updated selection/scroll
# position used while editing.
$ws.Columns("I").ColumnWidth = 9.75
$ws.Range("J2").Select()
